# update code tinh luong
$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "Đơn sale chính": add a new personal-customer order row, fix the
# discount rate/amount on the Filler order, and refresh the totals row.
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Row 8 (Trần Ái Mỹ / Tiêm Filler): source changed from "Khách cũ" to
# "Cá nhân" and the commission rate bumped from 10% to 15%.
$ws1.Range("F8").Value = "Cá nhân"
$ws1.Range("M8").Value = 0.15
$ws1.Range("N8").Value = 3150000

# Insert a new data row at row 9 (pushes the old "Tổng" row down to row 10).
$ws1.Range("A9").EntireRow.Insert()

$ws1.Range("A9").Value = "HD-LUXURY"
$ws1.Range("B9").Value = 596
$ws1.Range("C9").NumberFormat = "@"
$ws1.Range("C9").Value = "07-28-2024"
$ws1.Range("D9").Value = "SÓC TRĂNG"
$ws1.Range("E9").Value = "bảo trân"
$ws1.Range("F9").Value = "Cá nhân"
$ws1.Range("G9").Value = "Cắt mí"
$ws1.Range("H9").Value = 2000000
$ws1.Range("I9").Value = 0
$ws1.Range("J9").Value = 0
$ws1.Range("K9").Value = 2000000
$ws1.Range("L9").Value = 2000000
$ws1.Range("M9").Value = 0.1
$ws1.Range("N9").Value = 200000

# Refresh the "Tổng" row, now shifted to row 10.
$ws1.Range("B10").Value = 8
$ws1.Range("H10").Value = 120000000
$ws1.Range("K10").Value = 120000000
$ws1.Range("L10").Value = 117000000
$ws1.Range("N10").Value = 16570000

# ----------------------------------------------------------------------
# Sheet "Lương": the per-location breakdown gains explicit "Tổng công" /
# "Phụ cấp" / "Lương công tác" lines (previously only a single combined
# "Ngày công" / "Phụ cấp" pair existed for CẦN THƠ, and LONG XUYÊN /
# SÓC TRĂNG had none), and several totals were recalculated.
# ----------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

# Rows 2-3 become CẦN THƠ-specific labels with updated values.
$ws4.Range("A2").Value = "Tổng công tại CẦN THƠ"
$ws4.Range("B2").Value = 1
$ws4.Range("A3").Value = "Phụ cấp tại CẦN THƠ"
$ws4.Range("B3").Value = 35000

# Insert a new row 4: "Lương công tác tại CẦN THƠ".
$ws4.Range("A4").EntireRow.Insert()
$ws4.Range("A4").Value = "Lương công tác tại CẦN THƠ"
$ws4.Range("B4").Value = 350357.1428571428

# "Chiết khấu sale chính tại CẦN THƠ" shifted to row 6, value updated.
$ws4.Range("B6").Value = 10650000

# Insert 3 new rows at 14-16 for LONG XUYÊN "Tổng công / Phụ cấp / Lương công tác".
$ws4.Range("A14:A16").EntireRow.Insert()
$ws4.Range("A14").Value = "Tổng công tại LONG XUYÊN"
$ws4.Range("B14").Value = 0
$ws4.Range("A15").Value = "Phụ cấp tại LONG XUYÊN"
$ws4.Range("B15").Value = 0
$ws4.Range("A16").Value = "Lương công tác tại LONG XUYÊN"
$ws4.Range("B16").Value = 0

# Insert 2 new rows at 26-27 for SÓC TRĂNG "Tổng công / Phụ cấp".
$ws4.Range("A26:A27").EntireRow.Insert()
$ws4.Range("A26").Value = "Tổng công tại SÓC TRĂNG"
$ws4.Range("B26").Value = 21
$ws4.Range("A27").Value = "Phụ cấp tại SÓC TRĂNG"
$ws4.Range("B27").Value = 735000

# Remaining value-only updates at their final (post-insert) row numbers.
$ws4.Range("B28").Value = 7357500            # Lương cơ bản tại SÓC TRĂNG
$ws4.Range("B29").Value = 5920000            # Chiết khấu sale chính tại SÓC TRĂNG
$ws4.Range("B36").Value = -10000000          # Ứng lương tại SÓC TRĂNG
$ws4.Range("B37").Value = 11035357.14285714  # Tổng lương tại CẦN THƠ
$ws4.Range("B39").Value = 4632500            # Tổng lương tại SÓC TRĂNG
$ws4.Range("B40").Value = 15667857.14285714  # Tổng lương
